# gamma_ammo_weapon_table.xlsx
# ".338 federal as mutant hunter + 7.62x54r PP increased damage + text"
#
# 1) .338 Federal (row 3) re-tuned as a "mutant hunter" load: lower pen (G),
#    much higher raw damage (H) -> J/K recalc automatically via the existing
#    formulas.
# 2) 7.62x54r PP (row 36) gets a damage bump (H36) -> dependents recalc.
# 3) New ammo added: ammo_pkm_100 (row 43, AP type) with its own stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- 1) .338 Federal (row 3) : lower penetration, big damage boost ---------
$ws.Range("G3").Value = 0.18
$ws.Range("H3").Value = 3.78

# --- 2) 7.62x54r PP (row 36) : damage increase ------------------------------
$ws.Range("H36").Value = 1.73

# --- 3) New row: ammo_pkm_100 ----------------------------------------------
$ws.Range("A43").Value = "ammo_pkm_100"
$ws.Range("B43").Value = "AP"
$ws.Range("C43").Value = 5000
$ws.Range("D43").Formula = '=C43/30'
$ws.Range("E43").Formula = '=K43/D43'
$ws.Range("F43").Formula = '=G43/D43*100'
$ws.Range("G43").Value = 0.37
$ws.Range("H43").Value = 1.02
$ws.Range("I43").Value = 0.77
$ws.Range("J43").Formula = '=I43*H43'
$ws.Range("K43").Formula = '=J43*Feuil2!$B$1'

# Highlight the freshly-added row like a diff insertion (teal ammo id /
# pale-green price, both in a monospace face, vertically centred).
$ws.Range("A43").Font.Name = "Consolas"
$ws.Range("A43").Font.Family = 3
$ws.Range("A43").Font.Color = 11585870
$ws.Range("A43").VerticalAlignment = -4108

$ws.Range("C43").Font.Name = "Consolas"
$ws.Range("C43").Font.Family = 3
$ws.Range("C43").Font.Color = 11062965
$ws.Range("C43").VerticalAlignment = -4108

# Move the active selection off the old H37 cell onto H3, matching where the
# author was last working.
$ws.Range("H3").Select()
